# AFDP-2522: Update of foia drools rule files
#
# The "Save Case File Rules" rule table (Sheet1) has four rows whose
# CONDITION expressions (column C) use Groovy/SpEL-style null-safe
# navigation (queue?.name). These are simplified to plain property
# access (queue.name), and the two "Nullify" rules drop their now
# redundant "<field>EnterDate != null" clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28 - Set Biiling Enter Date
$ws.Range("C28").Value = "queue.name == 'Billing' && billingEnterDate == null"

# Row 29 - Nullify Billing Enter Date
$ws.Range("C29").Value = "queue.name != 'Billing'"

# Row 30 - Set Hold Enter Date
$ws.Range("C30").Value = "queue.name == 'Hold' && holdEnterDate == null"

# Row 31 - Nullify Hold Enter Date
$ws.Range("C31").Value = "queue.name != 'Hold'"

# Update the active selection left on the sheet after editing.
$ws.Range("D31").Select() | Out-Null
